$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: split the old combined "Diff test: ..." string into a short label (A3)
# plus the actual command in a new cell (B3).
$ws.Range("A3").Value = "Diff test:"
$ws.Range("B3").Value = "xltablediff.py  --key ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"
$ws.Range("B3").Style = $ws.Range("A3").Style

# Row 4: same treatment for the merge test line.
$ws.Range("A4").Value = "Merge test:"
$ws.Range("B4").Value = "xltablediff.py  --key ID --merge Color test1old.xlsx test1new.xlsx --out test1merge.xlsx"
$ws.Range("B4").Style = $ws.Range("A4").Style

# Row 5: same treatment for the append test line.
$ws.Range("A5").Value = "Append test:"
$ws.Range("B5").Value = "xltablediff.py  --key ID --append test1old.xlsx test1new.xlsx --out test1append.xlsx"
$ws.Range("B5").Style = $ws.Range("A5").Style

# The selection moved from A2 to A3.
$ws.Range("A3").Select() | Out-Null
